# Adjusted testdata for TC verifyGetStartedPages
$wb = $excel.ActiveWorkbook

# "Testdata" is the first sheet (sheet1.xml / rId1)
$ws = $wb.Worksheets.Item("Testdata")
$ws.Activate()

# Row 2, column A held "loginUser1" - rename to the new TC name.
# (This also removes "loginUser1" from the shared-strings table and
# appends "verifyGetStartedPages", shifting the other string indices
# down by one - matching the rest of the diff automatically.)
$ws.Range("A2").Value = "verifyGetStartedPages"

# Selection on the Testdata sheet moved from C14 to B2.
$ws.Range("B2").Select()
